$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-05 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.494707487409747
$ws.Range("E2").Value = 0.009284006189337557

$ws.Range("D3").Value = 0.2455238719580063
$ws.Range("E3").Value = 0.01983045716015752

$ws.Range("D4").Value = 0.09921518746102923
$ws.Range("E4").Value = 0.005397941250313787

$ws.Range("D5").Value = 0.1019813462209867
$ws.Range("E5").Value = 0.006978772899098651

$ws.Range("D6").Value = 0.03020593174198809
$ws.Range("E6").Value = 0.002258001178087499

$ws.Range("D7").Value = 0.02836617520824281
$ws.Range("E7").Value = 0.004369825206991695

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = 0.01090114066486647
